$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text ("@") number format on price cells whose new values would otherwise
# be auto-parsed by Excel as numbers, so they stay as plain text strings (matching
# the source data which stores all prices as text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "41.206.44"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "2.248.63"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "302.53"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").Value = "90.97"
$ws.Range("E6").Value = "  +5.00%  "
$ws.Range("D7").Value = "0.519"
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("D10").Value = "54.13"
$ws.Range("E10").Value = "  +10.34%  "
$ws.Range("D11").Value = "31.79"
$ws.Range("E11").Value = "  +7.09%  "
$ws.Range("D12").Value = "0.0794"
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("E13").Value = "  +3.45%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "2.594.46"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "14.10"
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.262.69"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").Value = "41.115.54"
$ws.Range("E19").Value = "  +3.33%  "
$ws.Range("D20").Value = "11.87"
$ws.Range("E20").Value = "  +5.50%  "
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("E22").Value = "  +1.68%  "
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").Value = "240.59"
$ws.Range("E24").Value = "  +2.09%  "
$ws.Range("E25").Value = "  +4.10%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +2.87%  "
$ws.Range("D28").Value = "23.74"
$ws.Range("E28").Value = "  +6.10%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  +1.51%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "9.57"
$ws.Range("E30").Value = "  +4.81%  "
$ws.Range("D31").Value = "158.94"
$ws.Range("E31").Value = "  +2.55%  "
$ws.Range("D32").Value = "33.06"
$ws.Range("E32").Value = "  +4.71%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  +6.19%  "
$ws.Range("D35").Value = "0.0730"
$ws.Range("E35").Value = "  +2.99%  "
$ws.Range("E36").Value = "  +7.80%  "
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "16.50"
$ws.Range("E39").Value = "  +6.96%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.103"
$ws.Range("E40").Value = "  +5.68%  "
$ws.Range("D41").Value = "1.77"
$ws.Range("E41").Value = "  +7.12%  "
$ws.Range("E42").Value = "  +5.32%  "
$ws.Range("D43").Value = "2.075.49"
$ws.Range("E43").Value = "  -2.16%  "
$ws.Range("E44").Value = "  +13.21%  "
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("D46").Value = "10.25"
$ws.Range("E46").Value = "  +6.69%  "
$ws.Range("E47").Value = "  +11.12%  "
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").Value = "  -4.94%  "
$ws.Range("D49").Value = "2.468.54"
$ws.Range("E49").Value = "  +2.16%  "
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("E51").Value = "  +4.31%  "
